$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row of the worksheet.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# For every data row whose "is_new" flag (column E) is TRUE, the row only
# had columns A-F populated (with F = "N/A"). Fill the remaining image/type
# columns (G-N) with "N/A" as well, matching the other "new word" rows.
for ($r = 2; $r -le $lastRow; $r++) {
    $isNew = $ws.Cells.Item($r, 5).Value2
    if ($isNew -eq $true) {
        $ws.Range($ws.Cells.Item($r, 7), $ws.Cells.Item($r, 14)).Value = "N/A"
    }
}
